# "adding averages and more checks"
#
# Training Dashboard sheet:
#  - Header band (row 2) and the big title (row 1) get a bold white font
#    (they share the same bold font definition, so both are touched the
#    same way).
#  - Row 3 ("Endangered by Electricity A safety Training (SOPs)"):
#      PERIOD TO EXPIRE  -43  -> -51
#      LAST UPDATE       08-Sep-2025 -> 16-Sep-2025

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# --- Title + header row: bold white font ---------------------------------
$title = $ws.Range("A1")
$title.Font.Bold = $true
$title.Font.Size = 11
$title.Font.Color = 0xFFFFFF

$hdr = $ws.Range("A2:K2")
$hdr.Font.Bold = $true
$hdr.Font.Color = 0xFFFFFF

# --- Row 3 data refresh ----------------------------------------------------
$ws.Range("H3").Value = -51

# Force text so the date stays a literal label ("16-Sep-2025"), matching how
# the other date columns in this row (F3/G3) are stored, instead of Excel
# auto-converting the typed string into a date serial number.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "16-Sep-2025"
